# Retrieve population from database (reload) Progress
#
# 1. "Parameters" sheet: Generations (B3) 10 -> 250; active sheet/selection
#    moves back to "Project" (tabSelected) with selection at E10 on Parameters.
# 2. "Project" sheet: add a new row "continue_run" / "No"; becomes the
#    tab-selected sheet with selection at B5.

$wb = $excel.ActiveWorkbook

# --- Parameters sheet: Generations 10 -> 250 ---
$params = $wb.Worksheets.Item("Parameters")
$params.Range("B3").Value = 250

# --- Project sheet: append continue_run / No row ---
$project = $wb.Worksheets.Item("Project")
$project.Range("A6").Value = "continue_run"
$project.Range("B6").Value = "No"

# --- Selections / active sheet bookkeeping ---
# Parameters' new selection (no longer the active/tabSelected sheet)
$params.Range("E10").Select()

# Project becomes the active (tabSelected) sheet, selection on B5
$project.Activate()
$project.Range("B5").Select()
